$p = $ppt.ActivePresentation
$ds = $p.Designs
$ds | Get-Member
